$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp stored for the 2025-06-02 price check
$ws.Range("A10").Value = 45810.39396888889

# Append the new price check recorded on 2025-06-03
$ws.Range("A11").Value = 45811.3938219378
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C11").Value = "1Kg"
$ws.Range("D11").Value = "15,41€"
